$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("D1").Value = 4
$ws.Range("E1").Value = 5

# Row 2
$ws.Range("B2").Value = 1

# Row 3
$ws.Range("C3").Value = 1

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("D4").Value = 1

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("E5").Value = 1

# Update selection to match target (active cell E5)
$ws.Range("E5").Select()
